$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Ian%Leonard-Lorant%NULL%1,                  Xavier%Delabranche%NULL%3,                  Xavier%Delabranche%NULL%0,                  François%Severac%NULL%3,                  François%Severac%NULL%0,                  Julie%Helms%NULL%0,                  Julie%Helms%NULL%0,                  Coralie%Pauzet%NULL%1,                  Olivier%Collange%NULL%1,                  Françis%Schneider%NULL%1,                  Aissam%Labani%NULL%2,                  Aissam%Labani%NULL%0,                  Pascal%Bilbault%NULL%1,                  Sébastien%Moliere%NULL%2,                  Sébastien%Moliere%NULL%0,                  Pierre%Leyendecker%NULL%1,                  Catherine%Roy%NULL%2,                  Catherine%Roy%NULL%0,                  Mickaël%Ohana%mickael.ohana@gmail.com%3,                  Mickaël%Ohana%mickael.ohana@gmail.com%0]'
$ws.Range("E3").Value = '[Franck%Grillet%f1grillet@chu-besancon.fr%1,                  Julien%Behr%NULL%2,                  Julien%Behr%NULL%0,                  Paul%Calame%NULL%1,                  Sébastien%Aubry%NULL%1,                  Eric%Delabrousse%NULL%2,                  Eric%Delabrousse%NULL%0]'
$ws.Range("E4").Value = '[Alban%Gervaise%alban.gervaise@hotmail.fr%1,                  Caroline%Bouzad%NULL%2,                  Caroline%Bouzad%NULL%0,                  Evelyne%Peroux%NULL%1,                  Carole%Helissey%NULL%1]'
$ws.Range("E5").Value = '[Hanny%Al-Samkari%NULL%1,                  Rebecca S.%Karp Leaf%NULL%1,                  Walter H.%Dzik%NULL%1,                  Jonathan C.T.%Carlson%NULL%1,                  Annemarie E.%Fogerty%NULL%1,                  Anem%Waheed%NULL%1,                  Katayoon%Goodarzi%NULL%1,                  Pavan K.%Bendapudi%NULL%1,                  Larissa%Bornikova%NULL%1,                  Shruti%Gupta%NULL%1,                  David E.%Leaf%NULL%1,                  David J.%Kuter%NULL%1,                  Rachel P.%Rosovsky%NULL%1]'
$ws.Range("E6").Value = '[Bin%Ren%xref no email%1,           Feifei%Yan%xref no email%1,           Zhouming%Deng%xref no email%1,           Sheng%Zhang%xref no email%0,           Lingfei%Xiao%xref no email%1,           Meng%Wu%xref no email%1,           Lin%Cai%xref no email%1]'
$ws.Range("E7").Value = '[Li%Zhang%xref no email%4,           Xiaokai%Feng%xref no email%1,           Danqing%Zhang%xref no email%1,           Chunguo%Jiang%xref no email%1,           Heng%Mei%xref no email%1,           Jing%Wang%xref no email%1,           Cuihong%Zhang%xref no email%1,           Hong%Li%xref no email%1,           Xiaoling%Xia%xref no email%1,           Shuangshuang%Kong%xref no email%1,           Jia%Liao%xref no email%1,           Huijun%Jia%xref no email%1,           Xueqin%Pang%xref no email%1,           Yue%Song%xref no email%1,           Ying%Tian%xref no email%1,           Bin%Wang%xref no email%1,           Chun%Wu%xref no email%1,           Hongliang%Yuan%xref no email%1,           Yongxing%Zhang%xref no email%1,           Yuman%Li%xref no email%1,           Wei%Sun%xref no email%1,           Yanting%Zhang%xref no email%1,           Shuangshuang%Zhu%xref no email%1,           Shuyuan%Wang%xref no email%1,           Yuji%Xie%xref no email%1,           Shuping%Ge%xref no email%1,           Liming%Zhang%xref no email%1,           Yu%Hu%xref no email%0,           Mingxing%Xie%xref no email%1]'
$ws.Range("E8").Value = '[Jean‐François%Llitjos%NULL%1,                  Maxime%Leclerc%NULL%1,                  Camille%Chochois%NULL%1,                  Jean‐Michel%Monsallier%NULL%1,                  Michel%Ramakers%NULL%1,                  Malika%Auvray%NULL%1,                  Karim%Merouani%NULL%1]'
$ws.Range("E9").Value = '[Julie%Helms%NULL%1,                  Charles%Tacquard%NULL%1,                  François%Severac%NULL%0,                  Ian%Leonard-Lorant%NULL%1,                  Mickaël%Ohana%NULL%0,                  Xavier%Delabranche%NULL%0,                  Hamid%Merdji%NULL%1,                  Raphaël%Clere-Jehl%NULL%1,                  Malika%Schenck%NULL%1,                  Florence%Fagot Gandet%NULL%1,                  Samira%Fafi-Kremer%NULL%0,                  Vincent%Castelain%NULL%1,                  Francis%Schneider%NULL%0,                  Lélia%Grunebaum%NULL%1,                  Eduardo%Anglés-Cano%NULL%1,                  Laurent%Sattler%NULL%1,                  Paul-Michel%Mertes%NULL%0,                  Ferhat%Meziani%Ferhat.Meziani@chru-strasbourg.fr%1,                  NULL%NULL%NULL%0]'
$ws.Range("E10").Value = '[P.%Demelo-Rodríguez%NULL%1,                  E.%Cervilla-Muñoz%NULL%1,                  L.%Ordieres-Ortega%NULL%1,                  A.%Parra-Virto%NULL%1,                  M.%Toledano-Macías%NULL%1,                  N.%Toledo-Samaniego%NULL%1,                  A.%García-García%NULL%1,                  I.%García-Fernández-Bravo%NULL%1,                  Z.%Ji%NULL%1,                  J.%de-Miguel-Diez%NULL%1,                  L.A.%Álvarez-Sala-Walther%NULL%1,                  J.%del-Toro-Cervera%NULL%1,                  F.%Galeano-Valle%paco.galeano.valle@gmail.com%1]'
$ws.Range("E11").Value = '[Songping%Cui%NULL%1,                  Shuo%Chen%NULL%1,                  Xiunan%Li%NULL%1,                  Shi%Liu%NULL%1,                  Feng%Wang%NULL%0]'
$ws.Range("E12").Value = '[Saskia%Middeldorp%NULL%1,                  Michiel%Coppens%NULL%1,                  Thijs F.%van Haaps%NULL%1,                  Merijn%Foppen%NULL%1,                  Alexander P.%Vlaar%NULL%1,                  Marcella C.A.%Müller%NULL%1,                  Catherine C.S.%Bouman%NULL%1,                  Ludo F.M.%Beenen%NULL%1,                  Ruud S.%Kootte%NULL%1,                  Jarom%Heijmans%NULL%1,                  Loek P.%Smits%NULL%1,                  Peter I.%Bonta%NULL%1,                  Nick%van Es%NULL%1]'
$ws.Range("E13").Value = '[Julien%Poissy%xref no email%1,           Julien%Goutay%xref no email%1,           Morgan%Caplan%xref no email%1,           Erika%Parmentier%xref no email%1,           Thibault%Duburcq%xref no email%1,           Fanny%Lassalle%xref no email%1,           Emmanuelle%Jeanpierre%xref no email%1,           Antoine%Rauch%xref no email%1,           Julien%Labreuche%xref no email%1,           Sophie%Susen%xref no email%1,           Nicolas%Cousin%xref no email%1,           Arthur%Durand%xref no email%1,           Ahmed%El Kalioubie%xref no email%1,           Raphaël%Favory%xref no email%1,           Patrick%Girardie%xref no email%1,           Marion%Houard%xref no email%1,           Emmanuelle%Jaillette%xref no email%1,           Mercé%Jourdain%xref no email%1,           Geoffrey%Ledoux%xref no email%1,           Daniel%Mathieu%xref no email%1,           Anne-Sophie%Moreau%xref no email%1,           Christopher%Niles%xref no email%1,           Saad%Nseir%xref no email%1,           Thierry%Onimus%xref no email%1,           Sébastien%Préau%xref no email%1,           Laurent%Robriquet%xref no email%1,           Anahita%Rouzé%xref no email%1,           Arthur%Simonnet%xref no email%1,           Sophie%Six%xref no email%1,           Aurélia%Toussaint%xref no email%1,           Annabelle%Dupont%xref no email%1,           Anne%Bauters%xref no email%1,           Christophe%Zawadzki%xref no email%1,           Camille%Paris%xref no email%1,           Nathalie%Trillot%xref no email%1,           Bénédicte%Wibaut%xref no email%1,           Audrey%Hochart%xref no email%1,           Catherine%Marichez%xref no email%1,           Vincent%Dalibard%xref no email%1,           Sandrine%Vanderziepe%xref no email%1,           Laureline%Bourgeois%xref no email%1,           Anaïs%Gaul%xref no email%1,           Aurélie%Jospin%xref no email%1,           Nataliia%Stepina%xref no email%1,           Bénédicte%Pradines%xref no email%1,           Antoine%Tournoys%xref no email%1,           Thierry%Brousseau%xref no email%1,           Martine%Rémy%xref no email%1,           Antoine%Hutt%xref no email%1]'
$ws.Range("E14").Value = '[Florian%Bompard%NULL%1,                  Hippolyte%Monnier%NULL%1,                  Ines%Saab%NULL%1,                  Mickael%Tordjman%NULL%1,                  Hendy%Abdoul%NULL%1,                  Laure%Fournier%NULL%1,                  Olivier%Sanchez%NULL%1,                  Christine%Lorut%NULL%1,                  Guillaume%Chassagnon%NULL%1,                  Marie-Pierre%Revel%NULL%1]'
$ws.Range("E15").Value = '[Marco%Cattaneo%NULL%1,                  Elena M.%Bertinato%NULL%1,                  Simone%Birocchi%NULL%1,                  Carolina%Brizio%NULL%1,                  Daniele%Malavolta%NULL%1,                  Marco%Manzoni%NULL%1,                  Gesualdo%Muscarella%NULL%1,                  Michela%Orlandi%NULL%1]'
